{"js": "// Helper: wrap a fragment of <w:p> elements into a full OOXML \"flat package\"\n// document so it can be fed to Range.insertOoxml().\nfunction wrapBodyFragment(bodyXml) {\n  return (\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" ' +\n    'pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    '<pkg:xmlData>' +\n    '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n    '<w:body>' + bodyXml + '</w:body></w:document>' +\n    '</pkg:xmlData></pkg:part></pkg:package>'\n  );\n}\n\nconst body = context.document.body;\n\n// ---------------------------------------------------------------------\n// Change 1: \"Nesse exato momento Gabriel percebe nos olhos de L\u00facifer...\"\n// gets a new leading sentence (\"Ao olhar para L\u00facifer, \") as its own run,\n// with the rest of the original sentence kept (minus the old opener) in a\n// second run.\n// ---------------------------------------------------------------------\nconst oldSentence =\n  \"Nesse exato momento Gabriel percebe nos olhos de L\u00facifer que ele est\u00e1 \" +\n  \"tramando algo, mas ainda faltava ele descobrir o que era.\";\n\nconst hit1 = body.search(oldSentence, { matchCase: true });\nhit1.load(\"items\");\nawait context.sync();\n\nif (hit1.items.length === 0) {\n  throw new Error(\"Could not find the target sentence for change 1.\");\n}\n\nconst target1Paragraph = hit1.items[0].paragraphs.getFirst();\nconst target1Range = target1Paragraph.getRange(\"Whole\");\n\nconst ooxml1 = wrapBodyFragment(\n  '<w:p>' +\n    '<w:r><w:t xml:space=\"preserve\">Ao olhar para L\u00facifer, </w:t></w:r>' +\n    '<w:r><w:t>Gabriel percebe nos olhos de L\u00facifer que ele est\u00e1 tramando algo, ' +\n    'mas ainda faltava ele descobrir o que era.</w:t></w:r>' +\n  '</w:p>'\n);\n\ntarget1Range.insertOoxml(ooxml1, Word.InsertLocation.replace);\nawait context.sync();\n\n// ---------------------------------------------------------------------\n// Change 2: new dialogue is added right after \"...L\u00facifer disse enquanto\n// encarava todos os anjos\", the paragraph that used to carry the\n// \"_GoBack\" bookmark and the single-underline paragraph mark now becomes\n// an empty, still-underlined paragraph after the new content, and the\n// bookmark moves to the end of the newly added text.\n// ---------------------------------------------------------------------\nconst bookmarkAnchor =\n  \"-Manterei o tratado, para que os humanos possam se deliciar com toda \" +\n  \"aquela maldade, afinal \u00e9 o que me mantem forte. \u2013 L\u00facifer disse enquanto \" +\n  \"encarava todos os anjos\";\n\nconst hit2 = body.search(bookmarkAnchor, { matchCase: true });\nhit2.load(\"items\");\nawait context.sync();\n\nif (hit2.items.length === 0) {\n  throw new Error(\"Could not find the anchor paragraph for change 2.\");\n}\n\n// The paragraph right after the anchor is the one that currently just\n// holds the \"_GoBack\" bookmark plus the single-underline paragraph mark.\nconst anchorParagraph = hit2.items[0].paragraphs.getFirst();\nconst bookmarkParagraph = anchorParagraph.getNext();\nconst bookmarkRange = bookmarkParagraph.getRange(\"Whole\");\n\nconst ooxml2 = wrapBodyFragment(\n  '<w:p><w:r><w:t>E ent\u00e3o com um olhar sombrio em seu rosto L\u00facifer pergunta:</w:t></w:r></w:p>' +\n  '<w:p><w:r><w:t>-Algu\u00e9m aqui presente \u00e9 contra mim?</w:t></w:r></w:p>' +\n  '<w:p><w:r><w:lastRenderedPageBreak/><w:t>Todos os anjos se mantiveram quietos, exceto Gabriel:</w:t></w:r></w:p>' +\n  '<w:p><w:r><w:t>-L\u00facifer, n\u00e3o posso deixar com que voc\u00ea fa\u00e7a o que quiser com todos os anjos, ' +\n  'e muito menos com a humanidade- Gabriel disse colocando toda a sua indigna\u00e7\u00e3o e frustra\u00e7\u00e3o nessas palavras</w:t></w:r></w:p>' +\n  '<w:p><w:r><w:t>-Pois bem, j\u00e1 que voc\u00ea n\u00e3o gosta da maneira que o c\u00e9u \u00e9 e ser\u00e1 daqui para a frente, ' +\n  'ent\u00e3o voc\u00ea n\u00e3o precisa mais viver aqui. \u2013 L\u00facifer disse enquanto levantava suas m\u00e3os</w:t></w:r></w:p>' +\n  '<w:p><w:r><w:t>E como se fosse empurrado, Gabriel foi jogado do c\u00e9u, de forma que ele n\u00e3o teve como evitar, ' +\n  'e em seu \u00faltimo olhar consciente ele v\u00ea a porta dos c\u00e9us se fechando em brasas.</w:t></w:r>' +\n  '<w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/><w:bookmarkEnd w:id=\"0\"/></w:p>' +\n  '<w:p/>' +\n  '<w:p><w:pPr><w:rPr><w:u w:val=\"single\"/></w:rPr></w:pPr></w:p>'\n);\n\nbookmarkRange.insertOoxml(ooxml2, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Era 2047 / \"O tratado De Cinza\" - apply the commit's edits through the\n# Word COM object model.\n#\n# Strategy: locate each target paragraph by its (unique) literal text,\n# then replace that paragraph's Range (which in the COM model already\n# includes the trailing paragraph mark) with a small OOXML fragment via\n# Range.InsertXML(). This lets us control run boundaries, the\n# lastRenderedPageBreak marker and the \"_GoBack\" bookmark precisely,\n# the same way the diff shows them.\n\n$d = $word.ActiveDocument\n$wNs = \"xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'\"\n\nfunction Get-ParagraphIndexByText($doc, [string]$text) {\n    $paras = $doc.Paragraphs\n    for ($i = 1; $i -le $paras.Count; $i++) {\n        $p = $paras.Item($i)\n        $t = $p.Range.Text\n        # A paragraph's Range.Text keeps its trailing paragraph mark (CR,\n        # and CR+BEL for list items) - strip that before comparing.\n        $t = $t.TrimEnd([char]13, [char]7)\n        if ($t -eq $text) {\n            return $i\n        }\n    }\n    return -1\n}\n\n# ---------------------------------------------------------------------\n# Change 1: split the \"Nesse exato momento Gabriel percebe...\" paragraph\n# into two runs, with a new leading sentence.\n# ---------------------------------------------------------------------\n$oldSentence = \"Nesse exato momento Gabriel percebe nos olhos de L\u00facifer que ele est\u00e1 tramando algo, mas ainda faltava ele descobrir o que era.\"\n\n$idx1 = Get-ParagraphIndexByText $d $oldSentence\nif ($idx1 -lt 0) {\n    throw \"Could not locate the target paragraph for change 1.\"\n}\n\n$para1 = $d.Paragraphs.Item($idx1)\n$xml1 = @\"\n<w:p $wNs><w:r><w:t xml:space=\"preserve\">Ao olhar para L\u00facifer, </w:t></w:r><w:r><w:t>Gabriel percebe nos olhos de L\u00facifer que ele est\u00e1 tramando algo, mas ainda faltava ele descobrir o que era.</w:t></w:r></w:p>\n\"@\n$para1.Range.InsertXML($xml1) | Out-Null\n\n# ---------------------------------------------------------------------\n# Change 2: after \"...L\u00facifer disse enquanto encarava todos os anjos\",\n# add the new dialogue/ending paragraphs. The paragraph that used to hold\n# only the \"_GoBack\" bookmark (with a single-underline paragraph mark)\n# now becomes an empty paragraph right after a new blank line, and the\n# bookmark moves to the end of the newly added text.\n# ---------------------------------------------------------------------\n$anchorText = \"-Manterei o tratado, para que os humanos possam se deliciar com toda aquela maldade, afinal \u00e9 o que me mantem forte. \u2013 L\u00facifer disse enquanto encarava todos os anjos\"\n\n$anchorIdx = Get-ParagraphIndexByText $d $anchorText\nif ($anchorIdx -lt 0) {\n    throw \"Could not locate the anchor paragraph for change 2.\"\n}\n\n# The very next paragraph is the one that currently just carries the\n# \"_GoBack\" bookmark plus the single-underline paragraph mark.\n$bookmarkIdx = $anchorIdx + 1\n$bookmarkPara = $d.Paragraphs.Item($bookmarkIdx)\n\n$xml2 = @\"\n<w:p $wNs><w:r><w:t>E ent\u00e3o com um olhar sombrio em seu rosto L\u00facifer pergunta:</w:t></w:r></w:p>\n<w:p $wNs><w:r><w:t>-Algu\u00e9m aqui presente \u00e9 contra mim?</w:t></w:r></w:p>\n<w:p $wNs><w:r><w:lastRenderedPageBreak/><w:t>Todos os anjos se mantiveram quietos, exceto Gabriel:</w:t></w:r></w:p>\n<w:p $wNs><w:r><w:t>-L\u00facifer, n\u00e3o posso deixar com que voc\u00ea fa\u00e7a o que quiser com todos os anjos, e muito menos com a humanidade- Gabriel disse colocando toda a sua indigna\u00e7\u00e3o e frustra\u00e7\u00e3o nessas palavras</w:t></w:r></w:p>\n<w:p $wNs><w:r><w:t>-Pois bem, j\u00e1 que voc\u00ea n\u00e3o gosta da maneira que o c\u00e9u \u00e9 e ser\u00e1 daqui para a frente, ent\u00e3o voc\u00ea n\u00e3o precisa mais viver aqui. \u2013 L\u00facifer disse enquanto levantava suas m\u00e3os</w:t></w:r></w:p>\n<w:p $wNs><w:r><w:t>E como se fosse empurrado, Gabriel foi jogado do c\u00e9u, de forma que ele n\u00e3o teve como evitar, e em seu \u00faltimo olhar consciente ele v\u00ea a porta dos c\u00e9us se fechando em brasas.</w:t></w:r><w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/><w:bookmarkEnd w:id=\"0\"/></w:p>\n<w:p $wNs/>\n<w:p $wNs><w:pPr><w:rPr><w:u w:val=\"single\"/></w:rPr></w:pPr></w:p>\n\"@\n$bookmarkPara.Range.InsertXML($xml2) | Out-Null\n"}
